# Updated ligand/receptor TPM-derived expression metrics (Bmp6-Bmpr1a).
# Each tuple is (row, column, new value); row/column are 1-based sheet indices
# matching columns G..T (7..20) of the LR-pairs data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @(2, 7, 7.739652666666667),
    @(2, 8, 23.218958),
    @(2, 9, 0.6488398532974882),
    @(2, 10, 0.6488398532974882),
    @(2, 13, 6.101885666666667),
    @(2, 14, 18.305657),
    @(2, 15, 0.1093737608697887),
    @(2, 16, 0.1093737608697887),
    @(2, 17, 47.22647567171178),
    @(2, 18, 425.038281045406),
    @(2, 19, 0.07096605495734828),
    @(2, 20, 0.07096605495734828),
    @(3, 7, 7.739652666666667),
    @(3, 8, 23.218958),
    @(3, 9, 0.6488398532974882),
    @(3, 10, 0.6488398532974882),
    @(3, 14, 87.53628900000001),
    @(3, 15, 0.5230171820937495),
    @(3, 16, 0.5230171820937495),
    @(3, 17, 225.8334908629847),
    @(3, 18, 2032.501417766862),
    @(3, 19, 0.3393543917017741),
    @(3, 20, 0.3393543917017741),
    @(4, 7, 7.739652666666667),
    @(4, 8, 23.218958),
    @(4, 9, 0.6488398532974882),
    @(4, 10, 0.6488398532974882),
    @(4, 11, 2),
    @(4, 12, 0.6666666666666666),
    @(4, 13, 0.146644),
    @(4, 14, 0.439932),
    @(4, 15, 0.002628532664354407),
    @(4, 16, 0.002628532664354407),
    @(4, 17, 1.134973625650667),
    @(4, 18, 10.214762630856),
    @(4, 19, 0.001705496748327369),
    @(4, 20, 0.001705496748327369),
    @(5, 7, 7.739652666666667),
    @(5, 8, 23.218958),
    @(5, 9, 0.6488398532974882),
    @(5, 10, 0.6488398532974882),
    @(5, 13, 15.02284966666667),
    @(5, 14, 45.068549),
    @(5, 15, 0.2692783275177917),
    @(5, 16, 0.2692783275177917),
    @(5, 17, 116.2716384835491),
    @(5, 18, 1046.444746351942),
    @(5, 19, 0.1747185105228369),
    @(5, 20, 0.174718510522837),
    @(6, 7, 7.739652666666667),
    @(6, 8, 23.218958),
    @(6, 9, 0.6488398532974882),
    @(6, 10, 0.6488398532974882),
    @(6, 13, 5.288900666666667),
    @(6, 14, 15.866702),
    @(6, 15, 0.09480134312252211),
    @(6, 16, 0.09480134312252211),
    @(6, 17, 40.93425414850178),
    @(6, 18, 368.408287336516),
    @(6, 19, 0.06151088956402209),
    @(6, 20, 0.06151088956402209),
    @(7, 7, 7.739652666666667),
    @(7, 8, 23.218958),
    @(7, 9, 0.6488398532974882),
    @(7, 10, 0.6488398532974882),
    @(7, 13, 0.050258),
    @(7, 14, 0.150774),
    @(7, 15, 0.0009008537317934847),
    @(7, 16, 0.0009008537317934848),
    @(7, 17, 0.3889794637213334),
    @(7, 18, 3.500815173492),
    @(7, 19, 0.0005845098031793793),
    @(7, 20, 0.0005845098031793794),
    @(8, 9, 0.3053032463428815),
    @(8, 10, 0.3053032463428815),
    @(8, 13, 6.101885666666667),
    @(8, 14, 18.305657),
    @(8, 15, 0.1093737608697887),
    @(8, 16, 0.1093737608697887),
    @(8, 17, 22.22181060955267),
    @(8, 18, 199.996295485974),
    @(8, 19, 0.03339216425827652),
    @(8, 20, 0.03339216425827652),
    @(9, 9, 0.3053032463428815),
    @(9, 10, 0.3053032463428815),
    @(9, 14, 87.53628900000001),
    @(9, 15, 0.5230171820937495),
    @(9, 16, 0.5230171820937495),
    @(9, 18, 956.3673961873982),
    @(9, 19, 0.1596788435863277),
    @(9, 20, 0.1596788435863277),
    @(10, 9, 0.3053032463428815),
    @(10, 10, 0.3053032463428815),
    @(10, 11, 2),
    @(10, 12, 0.6666666666666666),
    @(10, 13, 0.146644),
    @(10, 14, 0.439932),
    @(10, 15, 0.002628532664354407),
    @(10, 16, 0.002628532664354407),
    @(10, 17, 0.5340472393360001),
    @(10, 18, 4.806425154024001),
    @(10, 19, 0.0008024995555457041),
    @(10, 20, 0.0008024995555457042),
    @(11, 9, 0.3053032463428815),
    @(11, 10, 0.3053032463428815),
    @(11, 13, 15.02284966666667),
    @(11, 14, 45.068549),
    @(11, 15, 0.2692783275177917),
    @(11, 16, 0.2692783275177917),
    @(11, 17, 54.71012377896867),
    @(11, 18, 492.391114010718),
    @(11, 19, 0.08221154756096348),
    @(11, 20, 0.08221154756096349),
    @(12, 9, 0.3053032463428815),
    @(12, 10, 0.3053032463428815),
    @(12, 13, 5.288900666666667),
    @(12, 14, 15.866702),
    @(12, 15, 0.09480134312252211),
    @(12, 16, 0.09480134312252211),
    @(12, 17, 19.26108671446267),
    @(12, 18, 173.349780430164),
    @(12, 19, 0.0289431578129714),
    @(12, 20, 0.0289431578129714),
    @(13, 9, 0.3053032463428815),
    @(13, 10, 0.3053032463428815),
    @(13, 13, 0.050258),
    @(13, 14, 0.150774),
    @(13, 15, 0.0009008537317934847),
    @(13, 16, 0.0009008537317934848),
    @(13, 17, 0.183029282852),
    @(13, 18, 1.647263545668),
    @(13, 19, 0.0002750335687966503),
    @(13, 20, 0.0002750335687966504),
    @(14, 7, 0.5470016666666667),
    @(14, 8, 1.641005),
    @(14, 9, 0.04585690035963046),
    @(14, 10, 0.04585690035963046),
    @(14, 13, 6.101885666666667),
    @(14, 14, 18.305657),
    @(14, 15, 0.1093737608697887),
    @(14, 16, 0.1093737608697887),
    @(14, 17, 3.337741629476112),
    @(14, 18, 30.039674665285),
    @(14, 19, 0.005015541654163951),
    @(14, 20, 0.005015541654163951),
    @(15, 7, 0.5470016666666667),
    @(15, 8, 1.641005),
    @(15, 9, 0.04585690035963046),
    @(15, 10, 0.04585690035963046),
    @(15, 14, 87.53628900000001),
    @(15, 15, 0.5230171820937495),
    @(15, 16, 0.5230171820937495),
    @(15, 17, 15.96083199227167),
    @(15, 18, 143.647487930445),
    @(15, 19, 0.02398394680564777),
    @(15, 20, 0.02398394680564777),
    @(16, 7, 0.5470016666666667),
    @(16, 8, 1.641005),
    @(16, 9, 0.04585690035963046),
    @(16, 10, 0.04585690035963046),
    @(16, 11, 2),
    @(16, 12, 0.6666666666666666),
    @(16, 13, 0.146644),
    @(16, 14, 0.439932),
    @(16, 15, 0.002628532664354407),
    @(16, 16, 0.002628532664354407),
    @(16, 17, 0.08021451240666667),
    @(16, 18, 0.7219306116600001),
    @(16, 19, 0.000120536360481334),
    @(16, 20, 0.000120536360481334),
    @(17, 7, 0.5470016666666667),
    @(17, 8, 1.641005),
    @(17, 9, 0.04585690035963046),
    @(17, 10, 0.04585690035963046),
    @(17, 13, 15.02284966666667),
    @(17, 14, 45.068549),
    @(17, 15, 0.2692783275177917),
    @(17, 16, 0.2692783275177917),
    @(17, 17, 8.217523805749444),
    @(17, 18, 73.957714251745),
    @(17, 19, 0.01234826943399131),
    @(17, 20, 0.01234826943399131),
    @(18, 7, 0.5470016666666667),
    @(18, 8, 1.641005),
    @(18, 9, 0.04585690035963046),
    @(18, 10, 0.04585690035963046),
    @(18, 13, 5.288900666666667),
    @(18, 14, 15.866702),
    @(18, 15, 0.09480134312252211),
    @(18, 16, 0.09480134312252211),
    @(18, 17, 2.893037479501111),
    @(18, 18, 26.03733731551),
    @(18, 19, 0.004347295745528636),
    @(18, 20, 0.004347295745528635),
    @(19, 7, 0.5470016666666667),
    @(19, 8, 1.641005),
    @(19, 9, 0.04585690035963046),
    @(19, 10, 0.04585690035963046),
    @(19, 13, 0.050258),
    @(19, 14, 0.150774),
    @(19, 15, 0.0009008537317934847),
    @(19, 16, 0.0009008537317934848),
    @(19, 17, 0.02749120976333333),
    @(19, 18, 0.24742088787),
    @(19, 19, 0.0000413103598174551),
    @(19, 20, 0.00004131035981745509)
)

foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}
